$d = $word.ActiveDocument

# 1. Update the title paragraph: style Title -> Heading 1, text -> "Others Claim Template"
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "Others Claim Template"
$p1.Style = "Heading 1"

# 2. Remove the intro paragraph ("This report contains claim details for ... ")
$d.Paragraphs.Item(2).Range.Delete()

# 3. Remove the "Detailed Claim Information:" paragraph (contains a manual line break)
$d.Paragraphs.Item(4).Range.Delete()

# 4. Remove the now-empty single-cell "Insured Name" summary table
$d.Tables.Item(1).Delete()

# 5. Drop the explicit "Table Grid" style from the remaining details table
$d.Tables.Item(1).Style = $null

# 6. Rename the "Claim Number" label to "Claim No"
$d.Content.Find.Execute("Claim Number", $true, $false, $false, $false, $false, $true, 1, $false, "Claim No", 2)
